$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (cols C..H) ---
# This engine stores the xml "width" as ColumnWidth + 5/6 (rounded to the
# nearest 1/6th "pixel" step). We back-solve the ColumnWidth to land on
# (or as near as achievable to) the target width, and setting ColumnWidth
# explicitly also clears any bestFit flag on that column, matching the diff.
$ws.Columns.Item(3).ColumnWidth = 10.666666666666666   # -> width 11.5
$ws.Columns.Item(4).ColumnWidth = 9.833333333333334    # -> width ~10.6640625
$ws.Columns.Item(5).ColumnWidth = 11.0                 # -> width ~11.83203125
$ws.Columns.Item(6).ColumnWidth = 11.5                 # -> width ~12.33203125
$ws.Columns.Item(7).ColumnWidth = 7.0                  # -> width ~7.83203125
$ws.Columns.Item(8).ColumnWidth = 7.5                  # -> width ~8.33203125

# --- Row 13: new value + "0.00" number format ---
$ws.Range("C13").NumberFormat = "0.00"
$ws.Range("C13").Value = 0.56024894514767931
$ws.Range("D13").Value = 0.189

# --- Row 14: new value + "0.00" number format ---
$ws.Range("C14").NumberFormat = "0.00"
$ws.Range("C14").Value = 1.1781399176954725

# --- Rows 18 & 19: swap C:F values between the two rows ---
# NOTE: the ".Value" getter misbehaves in this runtime (it surfaces a
# reflection placeholder string instead of the cell's value), so reads use
# ".Value2" instead; writes continue to use ".Value" which works correctly.
$r18c = $ws.Range("C18").Value2
$r18d = $ws.Range("D18").Value2
$r18e = $ws.Range("E18").Value2
$r18f = $ws.Range("F18").Value2

$r19c = $ws.Range("C19").Value2
$r19d = $ws.Range("D19").Value2
$r19e = $ws.Range("E19").Value2
$r19f = $ws.Range("F19").Value2

$ws.Range("C18").Value = $r19c
$ws.Range("D18").Value = $r19d
$ws.Range("E18").Value = $r19e
$ws.Range("F18").Value = $r19f

$ws.Range("C19").Value = $r18c
$ws.Range("D19").Value = $r18d
$ws.Range("E19").Value = $r18e
$ws.Range("F19").Value = $r18f

# --- Remove the vestigial empty row 23 (only contained an empty K23 cell) ---
$ws.Range("K23").ClearContents()
